$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segment names currently in column A (row 2..20), in order.
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

# Insert a new column before B; this shifts the old B:E data
# (PercActivations, PercSegmentAreas, RelativeCAMImportance,
# PercActivationsRescaled) to C:F, keeping their values/format intact.
$ws.Columns.Item(2).Insert()

# New header for the inserted column, matching the style used by the
# other header cells (bold font, thin box border, centered/top aligned).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# Move the segment-name labels from column A into the new column B, and
# replace column A with a plain numeric row index (0-based).
for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $segments[$i]
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 1).Value = $i
}
